# Scheduled data refresh: update market-price / profit figures in the
# Moogle_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet holds static (non-formula) crafting-leve price data that the
# scheduled runner refreshes in place.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 11
$ws.Range("H11").Value = 84.583336
$ws.Range("K11").Value = 84.583336
$ws.Range("M11").Value = 55.416664
$ws.Range("I11").Value = 84.583336
# row 40
$ws.Range("J40").Value = 11971.571
$ws.Range("L40").Value = 11971.571
$ws.Range("H40").Value = 10565.111
$ws.Range("N40").Value = -12321.571
# row 51
$ws.Range("H51").Value = 6987
$ws.Range("L51").Value = 7198.8
$ws.Range("J51").Value = 7198.8
$ws.Range("N51").Value = -8166.8
# row 58
$ws.Range("K58").Value = 2379
$ws.Range("I58").Value = 793
$ws.Range("H58").Value = 5171.5
$ws.Range("M58").Value = -2229
$ws.Range("N58").Value = -28950
$ws.Range("J58").Value = 9550
$ws.Range("L58").Value = 28650
# row 74
$ws.Range("K74").Value = 19599.143
$ws.Range("I74").Value = 19599.143
$ws.Range("H74").Value = 18859.2
$ws.Range("M74").Value = -18663.143
# row 77
$ws.Range("H77").Value = 18859.2
$ws.Range("K77").Value = 97995.715
$ws.Range("I77").Value = 19599.143
$ws.Range("M77").Value = -93315.715
# row 137
$ws.Range("N137").Value = -37894.33199999999
$ws.Range("L137").Value = 32794.33199999999
$ws.Range("M137").Value = -1400.5002
$ws.Range("K137").Value = 3950.5002
$ws.Range("J137").Value = 10931.444
$ws.Range("H137").Value = 5437.381
$ws.Range("I137").Value = 1316.8334
# row 138
$ws.Range("N138").Value = -31137.9995
$ws.Range("J138").Value = 6952.6665
$ws.Range("L138").Value = 20857.9995
$ws.Range("H138").Value = 4572.2144

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1301.6552
$ws.Range("K2").Value = 1211.1305
$ws.Range("I2").Value = 1211.1305
$ws.Range("M2").Value = -1098.1305
# row 4
$ws.Range("N4").Value = -721.5
$ws.Range("H4").Value = 637.5
$ws.Range("J4").Value = 489.5
$ws.Range("L4").Value = 489.5
# row 44
$ws.Range("H44").Value = 48163
$ws.Range("J44").Value = 48163
$ws.Range("N44").Value = -49139
$ws.Range("L44").Value = 48163
# row 116
$ws.Range("H116").Value = 1301.6552
$ws.Range("K116").Value = 1211.1305
$ws.Range("I116").Value = 1211.1305
$ws.Range("M116").Value = 1082.8695
# row 119
$ws.Range("H119").Value = 183546.25
$ws.Range("J119").Value = 183546.25
$ws.Range("L119").Value = 183546.25
$ws.Range("N119").Value = -193222.25
# row 132
$ws.Range("J132").Value = 5658.1113
$ws.Range("M132").Value = -3107.272999999999
$ws.Range("H132").Value = 2976.2258
$ws.Range("L132").Value = 16974.3339
$ws.Range("K132").Value = 5637.272999999999
$ws.Range("I132").Value = 1879.091
$ws.Range("N132").Value = -22034.3339

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("I3").Value = 1211.1305
$ws.Range("K3").Value = 1211.1305
$ws.Range("H3").Value = 1301.6552
$ws.Range("M3").Value = -1097.1305
# row 5
$ws.Range("K5").Value = 0
$ws.Range("H5").Value = 75000
$ws.Range("M5").ClearContents()
$ws.Range("I5").Value = 0
# row 20
$ws.Range("M20").Value = -2528.3462
$ws.Range("H20").Value = 2923.6177
$ws.Range("K20").Value = 2775.3462
$ws.Range("J20").Value = 3405.5
$ws.Range("I20").Value = 2775.3462
$ws.Range("N20").Value = -3899.5
$ws.Range("L20").Value = 3405.5
# row 86
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -19080.666
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 16834.666
$ws.Range("J86").Value = 16834.666
$ws.Range("I86").Value = 0
$ws.Range("H86").Value = 16834.666
# row 89
$ws.Range("L89").Value = 84173.33
$ws.Range("H89").Value = 16834.666
$ws.Range("M89").ClearContents()
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 16834.666
$ws.Range("K89").Value = 0
$ws.Range("N89").Value = -95405.33
# row 99
$ws.Range("H99").Value = 3533.1667
$ws.Range("L99").Value = 4999.5
$ws.Range("N99").Value = -7995.5
$ws.Range("J99").Value = 4999.5
# row 105
$ws.Range("K105").Value = 4797.647
$ws.Range("L105").Value = 8715
$ws.Range("J105").Value = 8715
$ws.Range("I105").Value = 4797.647
$ws.Range("H105").Value = 6153.654
$ws.Range("M105").Value = -3050.647
$ws.Range("N105").Value = -12209
# row 132
$ws.Range("J132").Value = 99999
$ws.Range("N132").Value = -110119
$ws.Range("H132").Value = 99999
$ws.Range("L132").Value = 99999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("K22").Value = 608.0769
$ws.Range("I22").Value = 608.0769
$ws.Range("H22").Value = 2434.8635
$ws.Range("M22").Value = -258.0769
# row 107
$ws.Range("L107").Value = 2012.6
$ws.Range("J107").Value = 2012.6
$ws.Range("H107").Value = 798.8421
$ws.Range("N107").Value = -5852.6

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 40
$ws.Range("I40").Value = 624.6667
$ws.Range("J40").Value = 308.33334
$ws.Range("K40").Value = 2498.6668
$ws.Range("L40").Value = 1233.33336
$ws.Range("H40").Value = 585.125
$ws.Range("M40").Value = -2429.6668
$ws.Range("N40").Value = -1371.33336
# row 41
$ws.Range("J41").Value = 4205.75
$ws.Range("N41").Value = -13293.25
$ws.Range("H41").Value = 4205.75
$ws.Range("L41").Value = 12617.25
# row 97
$ws.Range("N97").Value = -14281.7
$ws.Range("J97").Value = 4429.9
$ws.Range("H97").Value = 4203.2
$ws.Range("L97").Value = 13289.7
# row 132
$ws.Range("M132").Value = -28207.25
$ws.Range("H132").Value = 3472.1333
$ws.Range("K132").Value = 30737.25
$ws.Range("I132").Value = 3415.25
# row 140
$ws.Range("H140").Value = 1520.0769
$ws.Range("N140").Value = -15846.0625
$ws.Range("J140").Value = 1828.6875
$ws.Range("L140").Value = 5486.0625

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 2797.3333
$ws.Range("K2").Value = 1897.125
$ws.Range("J2").Value = 9999
$ws.Range("L2").Value = 9999
$ws.Range("N2").Value = -10225
$ws.Range("I2").Value = 1897.125
$ws.Range("M2").Value = -1784.125
# row 113
$ws.Range("M113").Value = -451.3157000000001
$ws.Range("K113").Value = 2621.3157
$ws.Range("N113").Value = -8981.9
$ws.Range("H113").Value = 3318.0688
$ws.Range("J113").Value = 4641.9
$ws.Range("I113").Value = 2621.3157
$ws.Range("L113").Value = 4641.9

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("M7").Value = -6162.5835
$ws.Range("I7").Value = 6274.5835
$ws.Range("K7").Value = 6274.5835
$ws.Range("H7").Value = 7237.5
# row 22
$ws.Range("K22").Value = 1071.2
$ws.Range("I22").Value = 1071.2
$ws.Range("N22").Value = -1960.5
$ws.Range("J22").Value = 1370.5
$ws.Range("L22").Value = 1370.5
$ws.Range("H22").Value = 1170.9667
$ws.Range("M22").Value = -776.2
# row 27
$ws.Range("I27").Value = 1071.2
$ws.Range("N27").Value = -1584.5
$ws.Range("L27").Value = 1370.5
$ws.Range("M27").Value = -964.2
$ws.Range("H27").Value = 1170.9667
$ws.Range("J27").Value = 1370.5
$ws.Range("K27").Value = 1071.2
# row 46
$ws.Range("L46").Value = 5110.8887
$ws.Range("M46").Value = -1496.8889
$ws.Range("K46").Value = 1684.8889
$ws.Range("N46").Value = -5486.8887
$ws.Range("I46").Value = 1684.8889
$ws.Range("H46").Value = 3397.889
$ws.Range("J46").Value = 5110.8887
# row 55
$ws.Range("H55").Value = 1731.909
$ws.Range("I55").Value = 821.2857
$ws.Range("M55").Value = -648.2857
$ws.Range("K55").Value = 821.2857
# row 126
$ws.Range("I126").Value = 6274.5835
$ws.Range("H126").Value = 7237.5
$ws.Range("M126").Value = -16353.7505
$ws.Range("K126").Value = 18823.7505

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("I81").Value = 1030.1666
$ws.Range("H81").Value = 1030.1666
$ws.Range("M81").Value = -999.3332
$ws.Range("K81").Value = 2060.3332
# row 84
$ws.Range("I84").Value = 1030.1666
$ws.Range("H84").Value = 1030.1666
$ws.Range("K84").Value = 10301.666
$ws.Range("M84").Value = -4997.666000000001
# row 100
$ws.Range("N100").Value = -2971
$ws.Range("M100").Value = -279
$ws.Range("L100").Value = 1889
$ws.Range("K100").Value = 820
$ws.Range("J100").Value = 944.5
$ws.Range("H100").Value = 766.3333
$ws.Range("I100").Value = 410
# row 107
$ws.Range("K107").Value = 3387.4614
$ws.Range("L107").Value = 4614
$ws.Range("J107").Value = 1538
$ws.Range("M107").Value = -1467.4614
$ws.Range("H107").Value = 1225.3529
$ws.Range("I107").Value = 1129.1538
$ws.Range("N107").Value = -8454
# row 108
$ws.Range("J108").Value = 179883.33
$ws.Range("N108").Value = -187563.33
$ws.Range("L108").Value = 179883.33
$ws.Range("H108").Value = 179883.33
# row 126
$ws.Range("I126").Value = 2741.6
$ws.Range("H126").Value = 3236.0715
$ws.Range("M126").Value = -5754.799999999999
$ws.Range("K126").Value = 8224.799999999999
# row 132
$ws.Range("J132").Value = 5516.1875
$ws.Range("M132").Value = -5883.4547
$ws.Range("H132").Value = 3689.9387
$ws.Range("L132").Value = 16548.5625
$ws.Range("K132").Value = 8413.4547
$ws.Range("I132").Value = 2804.4849
$ws.Range("N132").Value = -21608.5625
